$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Expenditures by Team")
$ws.Activate()

# Update costs for the two RF TXRX modules (both the "by team" summary table AC:AG
# and the main part table B:I), and the corresponding Digikey part numbers.

# Team summary table (AC:AG) rows 10-11
$ws.Range("AE10").Value = 18.19
$ws.Range("AE11").Value = 18.19

# Main part table rows 19-20
$ws.Range("C19").Value = 18.19
$ws.Range("E19").Value = "602-1559-ND"

$ws.Range("C20").Value = 18.19
$ws.Range("E20").Value = "602-1560-ND"

# Order number update
$ws.Range("I22").Value = 57396497

# Update the sheet view (scroll position + selection) to match the committed state
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Range("AC10:AF11").Select()
